$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-08-27 10:45:30"

$zhcn.Range("H4").Value = "2016-08-27 10:45:26"
$zhcn.Range("K4").Value = "2016-08-27 10:45:44"

$dede.Range("K4").Value = "2016-08-27 10:45:51"
